$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 1092, shifting existing rows 1092:1163 down to 1095:1166
$ws.Rows.Item(1092).Insert()
$ws.Rows.Item(1093).Insert()
$ws.Rows.Item(1094).Insert()

# Populate new row 1092
$ws.Cells.Item(1092, 1).Value = 1
$ws.Cells.Item(1092, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(1092, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(1092, 4).Value = 44931
$ws.Cells.Item(1092, 5).Value = 15
$ws.Cells.Item(1092, 6).Value = 100112033
$ws.Cells.Item(1092, 7).Value = 'Lechuga'
$ws.Cells.Item(1092, 8).Value = 'Conconina(o)'
$ws.Cells.Item(1092, 9).Value = 'Segunda'
$ws.Cells.Item(1092, 10).Value = 1200
$ws.Cells.Item(1092, 11).Value = 400
$ws.Cells.Item(1092, 12).Value = 500
$ws.Cells.Item(1092, 13).Value = 450
$ws.Cells.Item(1092, 14).Value = '$/unidad'
$ws.Cells.Item(1092, 15).Value = 'Provincia de Chacabuco'
$ws.Cells.Item(1092, 16).Value = 450
$ws.Cells.Item(1092, 17).Value = 1
$ws.Cells.Item(1092, 18).Value = 'Hortaliza'

# Populate new row 1093
$ws.Cells.Item(1093, 1).Value = 1
$ws.Cells.Item(1093, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(1093, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(1093, 4).Value = 44931
$ws.Cells.Item(1093, 5).Value = 15
$ws.Cells.Item(1093, 6).Value = 100112033
$ws.Cells.Item(1093, 7).Value = 'Lechuga'
$ws.Cells.Item(1093, 8).Value = 'Escarola'
$ws.Cells.Item(1093, 9).Value = 'Primera'
$ws.Cells.Item(1093, 10).Value = 330
$ws.Cells.Item(1093, 11).Value = 3500
$ws.Cells.Item(1093, 12).Value = 4000
$ws.Cells.Item(1093, 13).Value = 3758
$ws.Cells.Item(1093, 14).Value = '$/caja 12 unidades'
$ws.Cells.Item(1093, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(1093, 16).Value = 313
$ws.Cells.Item(1093, 17).Value = 12
$ws.Cells.Item(1093, 18).Value = 'Hortaliza'

# Populate new row 1094
$ws.Cells.Item(1094, 1).Value = 1
$ws.Cells.Item(1094, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(1094, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(1094, 4).Value = 44931
$ws.Cells.Item(1094, 5).Value = 15
$ws.Cells.Item(1094, 6).Value = 100112033
$ws.Cells.Item(1094, 7).Value = 'Lechuga'
$ws.Cells.Item(1094, 8).Value = 'Escarola'
$ws.Cells.Item(1094, 9).Value = 'Segunda'
$ws.Cells.Item(1094, 10).Value = 200
$ws.Cells.Item(1094, 11).Value = 3500
$ws.Cells.Item(1094, 12).Value = 4000
$ws.Cells.Item(1094, 13).Value = 3750
$ws.Cells.Item(1094, 14).Value = '$/caja 18 unidades'
$ws.Cells.Item(1094, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(1094, 16).Value = 208
$ws.Cells.Item(1094, 17).Value = 18
$ws.Cells.Item(1094, 18).Value = 'Hortaliza'
